# Apply crypto price/volume updates per commit "Updated cryptos list on Fri Jun  7 10:15:14 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '71.303.95'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.61%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.810.05'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.93%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '702.85'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.45%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.74'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.40%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.808.87'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.93%  '

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.06%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.524'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.22%  '

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.29%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.50'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.87%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.470'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.75%  '

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.38%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.00'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.57%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.451.66'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.94%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.820.80'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.82%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '71.295.44'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.56%  '

# Row 18
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.17'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.19%  '

# Row 19
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.46'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.63%  '

# Row 20
$ws.Range("B20").Value = 'TRON'
$ws.Range("C20").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.115'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.19%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '514.87'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +4.85%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.64'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.31%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.716'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.05%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.06'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.19%  '

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.67%  '

# Row 26
$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.13'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.06%  '

# Row 27
$ws.Range("B27").Value = 'WrappedeETH'
$ws.Range("C27").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.957.94'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.11%  '

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.92%  '

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.11%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.03'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.73%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.03'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -4.70%  '

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.50%  '

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.62%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.06'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.03%  '

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.71%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.15'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.18%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.771.33'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.80%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.998'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.22%  '

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.16%  '

# Row 40
$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.28'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.73%  '

# Row 41
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.38'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.05%  '

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.82%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.28'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.11%  '

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.00%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '174.17'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +6.70%  '

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.06%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.000310'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.74%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '49.39'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.51%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '425.37'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.52%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.38'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.85%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.54'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.79%  '
